$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): B1 already holds 0 (generator id). Add C1 = 1
# (second generator id), carrying the same bold/centered/bordered header
# style that B1 already has. ---
$ws.Cells.Item(1, 3).Value = 1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats

# --- Column A (hour index) keeps the header-like style that was already on
# A2 in the original file; stretch that style down through A25 before
# overwriting the values for the full 24-hour profile. ---
$ws.Range("A2").Copy()
$ws.Range("A2:A25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# p_mw values (generator 0) per hour 0..23 and the constant max capacity
# (generator 1) column, rebuilt from the refreshed 24h CSV profile.
$bValues = @(0,0,0,0,0,0,0.1,0.2,0.25,0.33,0.41,0.65,0.82,1,0.9399999999999999,0.84,0.71,0.58,0.29,0.2,0.1,0,0,0)

for ($hour = 0; $hour -lt 24; $hour++) {
    $row = $hour + 2
    $ws.Cells.Item($row, 1).Value = $hour
    $ws.Cells.Item($row, 2).Value = $bValues[$hour]
    $ws.Cells.Item($row, 3).Value = 200
}
